$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84, pushing existing rows 84..122 down to 85..123
$ws.Rows.Item(84).Insert()

# Populate the new row 84 with the latest weekly entry
$ws.Cells.Item(84, 1).Value = 8
$ws.Cells.Item(84, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(84, 3).Value = "Coquimbo"
$ws.Cells.Item(84, 4).Value = 44489
$ws.Cells.Item(84, 5).Value = 4
$ws.Cells.Item(84, 6).Value = 100112037
$ws.Cells.Item(84, 7).Value = "Cebollín"
$ws.Cells.Item(84, 8).Value = "Sin especificar"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 3200
$ws.Cells.Item(84, 11).Value = 900
$ws.Cells.Item(84, 12).Value = 1000
$ws.Cells.Item(84, 13).Value = 950
$ws.Cells.Item(84, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(84, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(84, 16).Value = 158
$ws.Cells.Item(84, 17).Value = 6
$ws.Cells.Item(84, 18).Value = "Hortaliza"
